$wb = $excel.ActiveWorkbook
$wsUI = $wb.Worksheets.Item("UI")
$wsSTR = $wb.Worksheets.Item("STR")

# Update row 36, column C: "Nguyên liệu" -> "Tiêu hao"
# Order of writes below matters: it determines the order new strings are
# appended to the shared string table, which must match indices 612-617.
$wsUI.Range("C61").Value = "Gỡ"
$wsUI.Range("B61").Value = "UnEquip"
$wsUI.Range("A61").Value = "UI_UNEQUIP"

$wsUI.Range("C36").Value = "Tiêu hao"

$wsUI.Range("A62").Value = "UI_ENCHANCE"
$wsUI.Range("B62").Value = "Enchance"
$wsUI.Range("C62").Value = "Tăng cấp"

# Restore the selection on the UI sheet before switching the active tab,
# so the saved selection for sheet1 reflects the new cell.
$wsUI.Range("C57").Select()

# Move the topLeftCell scroll position for the UI sheet.
$win = $excel.ActiveWindow
$win.ScrollRow = 44
$win.ScrollColumn = 1

# Make the STR sheet the active tab (was previously UI).
$wsSTR.Activate()

Write-Host "done"
